$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 2 (shifts existing rows 2-11 down to 3-12) ---
$ws.Rows.Item(2).Insert()

# Fill the new row 2 (Specification of RTuinOS API)
$ws.Range("A2").Value = 41164
$ws.Range("C2").Value = 2.5
$ws.Range("D2").Value = "Specification of RTuinOS API"

# --- Update row 9's description (text revised) ---
$ws.Range("D9").Value = "Continuation of implementation rtos.c. Implementation problems with first suspend operation; New implementation concept made but not proven or implemented yet"

# --- Append new rows 13-15 at the end ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = 41179
$ws.Range("B13").Value = 3
$ws.Range("D13").Value = "Debugging of current implementation. Still no success with two tasks plus idle"

$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = 41180
$ws.Range("B14").Value = 3.75
$ws.Range("C14").Value = 0.25
$ws.Range("D14").Value = "First success with two tasks plus idle but not yet proven by profound testing"

$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = 41183
$ws.Range("B15").Value = 2
$ws.Range("D15").Value = "First non-trivial test cases are running well"

# --- Update selection to match authored state ---
$ws.Range("E15").Select()
